$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the source data. It sorts before the
# existing row 514, so insert a fresh row there (pushing rows 514:569 down to
# 515:570) and populate it with the new observation.
$ws.Rows.Item(514).Insert()

$ws.Range("A514").Value = 10
$ws.Range("B514").Value = 'Vega Modelo de Temuco'
$ws.Range("C514").Value = 'La Araucanía'
$ws.Range("D514").Value = 45194
$ws.Range("E514").Value = 9
$ws.Range("F514").Value = 100112009
$ws.Range("G514").Value = 'Acelga'
$ws.Range("H514").Value = 'Sin especificar'
$ws.Range("I514").Value = 'Primera'
$ws.Range("J514").Value = 50
$ws.Range("K514").Value = 8000
$ws.Range("L514").Value = 8000
$ws.Range("M514").Value = 8000
$ws.Range("N514").Value = '$/docena de atados (12 kilos)'
$ws.Range("O514").Value = 'Provincia de Cautín'
$ws.Range("P514").Value = 667
$ws.Range("Q514").Value = 12
$ws.Range("R514").Value = 'Hortaliza'
